$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting existing rows 29-124 down to 30-125
$ws.Rows.Item(29).Insert()

$ws.Range("A29").Value = 5
$ws.Range("B29").Value = "Macroferia Regional de Talca"
$ws.Range("C29").Value = "Maule"
$ws.Range("D29").Value = 45177
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = 100112040
$ws.Range("G29").Value = "Cilantro"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 8000
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = 8000
$ws.Range("N29").Value = "$/caja 36 atados"
$ws.Range("O29").Value = "Región Metropolitana"
$ws.Range("P29").Value = 222
$ws.Range("Q29").Value = 36
$ws.Range("R29").Value = "Hortaliza"
